$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting for numeric-looking price strings
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '35.478.65'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.895.70'
$ws.Range("E3").Value = '  -0.79%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.80%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '247.38'
$ws.Range("E5").Value = '  -2.85%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '0.692'
$ws.Range("E6").Value = '  -3.84%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.93%  '

$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").Value = '44.09'
$ws.Range("E8").Value = '  +8.71%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.352'
$ws.Range("E9").Value = '  -4.10%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.0737'
$ws.Range("E10").Value = '  -3.82%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.0971'
$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '13.20'
$ws.Range("E12").Value = '  +3.51%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.171.23'
$ws.Range("E13").Value = '  -0.72%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.726'
$ws.Range("E14").Value = '  +0.84%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '4.93'
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.856.71'
$ws.Range("E16").Value = '  -2.97%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '35.445.67'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '73.36'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0822'
$ws.Range("E19").Value = '  -3.45%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '246.73'
$ws.Range("E20").Value = '  +1.26%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '12.86'
$ws.Range("E21").Value = '  -1.92%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.96'
$ws.Range("E22").Value = '  -2.79%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.86%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '2.54'
$ws.Range("E24").Value = '  +6.34%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.19'
$ws.Range("E25").Value = '  -10.88%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '165.47'
$ws.Range("E26").Value = '  -1.00%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '8.49'
$ws.Range("E27").Value = '  -2.38%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.36'
$ws.Range("E28").Value = '  -2.13%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.128'
$ws.Range("E29").Value = '  -3.87%  '

$ws.Range("B30").Value = 'EURNeutrino'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn'
$ws.Range("D30").Value = '4.128.49'
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("B31").Value = 'TrustWalletToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D31").Value = '1.82'
$ws.Range("E31").Value = '  +10.34%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.27'
$ws.Range("E32").Value = '  -2.20%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0582'
$ws.Range("E33").Value = '  -0.57%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.26'
$ws.Range("E34").Value = '  +1.10%  '

$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.86%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.851'
$ws.Range("E36").Value = '  -6.56%  '

$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.02'
$ws.Range("E37").Value = '  -1.20%  '

$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '1.59'
$ws.Range("E38").Value = '  -20.39%  '

$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = '17.26'
$ws.Range("E39").Value = '  +0.93%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '97.78'
$ws.Range("E40").Value = '  +0.96%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.0673'
$ws.Range("E41").Value = '  +3.67%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0213'
$ws.Range("E42").Value = '  -1.90%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '1.10'
$ws.Range("E43").Value = '  -1.94%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.289.85'
$ws.Range("E44").Value = '  -3.52%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '2.36'
$ws.Range("E45").Value = '  -2.07%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.0803'
$ws.Range("E46").Value = '  +7.55%  '

$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").Value = '2.40'
$ws.Range("E47").Value = '  -1.03%  '

$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("B49").Value = 'Gas'
$ws.Range("C49").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D49").Value = '12.05'
$ws.Range("E49").Value = '  +1.61%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '6.41'
$ws.Range("E50").Value = '  -5.14%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '43.19'
$ws.Range("E51").Value = '  -4.89%  '
